# Recompute/update the "adjacency" distance column (E) and a couple of
# lat/long source values (B13, C13) with refreshed figures, then leave the
# selection where the author left off (cell L30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E (adjacency distance) value refresh -------------------------
$ws.Range("E2").Value = 113.9789267
$ws.Range("E3").Value = 604.71490226709398
$ws.Range("E4").Value = 130.46883645572231
$ws.Range("E5").Value = 80.919601483103435
$ws.Range("E7").Value = 118.88249381413505
$ws.Range("E8").Value = 118.88249381413505
$ws.Range("E9").Value = 118.88249381413505
$ws.Range("E10").Value = 118.88249381413505
$ws.Range("E11").Value = 118.88249381413505
$ws.Range("E12").Value = 118.88249381413505

# --- Row 13: refreshed lat/long + distance --------------------------------
$ws.Range("B13").Value = 51.922169573927299
$ws.Range("C13").Value = 4.4098035896042003
$ws.Range("E13").Value = 83.802219927778751

$ws.Range("E14").Value = 77.959626151355195
$ws.Range("E15").Value = 80.664968396488433
$ws.Range("E16").Value = 123.43473173813405
$ws.Range("E17").Value = 136.44680478452347

# --- Leave the selection where the author finished up ---------------------
$ws.Range("L30").Select()
